$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BurndownChart")

$ws.Range("B1").Value = 5
$ws.Range("B14").Value = 0.2
$ws.Range("B16").Value = 0.3
$ws.Range("B17").Value = 0.5
$ws.Range("C17").Value = 0.5
$ws.Range("B18").Value = 0.5
$ws.Range("C18").Value = 0.5

$ws.Range("B20").Select()
